$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1989.5588
$ws.Range("I15").Value = 1989.5588
$ws.Range("K15").Value = 5968.6764
$ws.Range("M15").Value = -5799.6764

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H40").Value = 4500
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350

$ws.Range("H64").Value = 63928
$ws.Range("I64").Value = 85599.2
$ws.Range("K64").Value = 85599.2
$ws.Range("M64").Value = -85351.2

$ws.Range("H67").Value = 63928
$ws.Range("I67").Value = 85599.2
$ws.Range("K67").Value = 85599.2
$ws.Range("M67").Value = -84741.2

$ws.Range("H100").Value = 17964372
$ws.Range("I100").Value = 41834016
$ws.Range("J100").Value = 62136.625
$ws.Range("K100").Value = 41834016
$ws.Range("L100").Value = 62136.625
$ws.Range("M100").Value = -41833475
$ws.Range("N100").Value = -63218.625

$ws.Range("H113").Value = 7883.5713
$ws.Range("I113").Value = 7871.375
$ws.Range("J113").Value = 7899.8335
$ws.Range("K113").Value = 7871.375
$ws.Range("L113").Value = 7899.8335
$ws.Range("M113").Value = -4617.375
$ws.Range("N113").Value = -14407.8335

$ws.Range("H129").Value = 1688.4
$ws.Range("I129").Value = 1183.7142
$ws.Range("J129").Value = 2866
$ws.Range("K129").Value = 3551.1426
$ws.Range("L129").Value = 8598
$ws.Range("M129").Value = 1448.8574
$ws.Range("N129").Value = -18598

$ws.Range("H132").Value = 3065
$ws.Range("I132").Value = 3277.606
$ws.Range("K132").Value = 9832.818000000001
$ws.Range("M132").Value = -7302.818000000001

$ws.Range("H141").Value = 4362.421
$ws.Range("I141").Value = 3368.8276
$ws.Range("K141").Value = 10106.4828
$ws.Range("M141").Value = -4926.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 12286
$ws.Range("I39").Value = 14508
$ws.Range("J39").Value = 11175
$ws.Range("K39").Value = 14508
$ws.Range("L39").Value = 11175
$ws.Range("M39").Value = -13988
$ws.Range("N39").Value = -12215

$ws.Range("H40").Value = 18014
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6352

$ws.Range("H92").Value = 366683330
$ws.Range("J92").Value = 366683330
$ws.Range("L92").Value = 366683330
$ws.Range("N92").Value = -366688322

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2079.9565
$ws.Range("I20").Value = 720.2
$ws.Range("J20").Value = 3125.923
$ws.Range("K20").Value = 720.2
$ws.Range("L20").Value = 3125.923
$ws.Range("M20").Value = -473.2
$ws.Range("N20").Value = -3619.923

$ws.Range("H38").Value = 8599.75
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 8599.75
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 8599.75
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -9431.75

$ws.Range("H43").Value = 283087.25
$ws.Range("J43").Value = 283087.25
$ws.Range("L43").Value = 283087.25
$ws.Range("N43").Value = -283449.25

$ws.Range("H105").Value = 4357.12
$ws.Range("I105").Value = 3526.75
$ws.Range("K105").Value = 3526.75
$ws.Range("M105").Value = -1779.75

$ws.Range("H134").Value = 6282.8213
$ws.Range("I134").Value = 6311.1567
$ws.Range("J134").Value = 5993.8
$ws.Range("K134").Value = 18933.4701
$ws.Range("L134").Value = 17981.4
$ws.Range("M134").Value = -16398.4701
$ws.Range("N134").Value = -23051.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3859.4
$ws.Range("I31").Value = 2758.182
$ws.Range("K31").Value = 2758.182
$ws.Range("M31").Value = -2463.182

$ws.Range("H34").Value = 3859.4
$ws.Range("I34").Value = 2758.182
$ws.Range("K34").Value = 2758.182
$ws.Range("M34").Value = -2556.182

$ws.Range("H58").Value = 3424.6538
$ws.Range("J58").Value = 3635.25
$ws.Range("L58").Value = 3635.25
$ws.Range("N58").Value = -4041.25

$ws.Range("H107").Value = 10811.272
$ws.Range("I107").Value = 15149.866
$ws.Range("K107").Value = 15149.866
$ws.Range("M107").Value = -13229.866

$ws.Range("H132").Value = 12765.658
$ws.Range("I132").Value = 1520.2069
$ws.Range("K132").Value = 4560.620699999999
$ws.Range("M132").Value = -2030.620699999999

$ws.Range("H136").Value = 3424.6538
$ws.Range("J136").Value = 3635.25
$ws.Range("L136").Value = 10905.75
$ws.Range("N136").Value = -16005.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 6662.5713
$ws.Range("J124").Value = 9916.25
$ws.Range("L124").Value = 29748.75
$ws.Range("N124").Value = -39568.75

$ws.Range("H126").Value = 6693.6665
$ws.Range("I126").Value = 2472.7144
$ws.Range("J126").Value = 12603
$ws.Range("K126").Value = 7418.1432
$ws.Range("L126").Value = 37809
$ws.Range("M126").Value = -2478.1432
$ws.Range("N126").Value = -47689

$ws.Range("H129").Value = 9092615
$ws.Range("J129").Value = 12501825
$ws.Range("L129").Value = 37505475
$ws.Range("N129").Value = -37515475

$ws.Range("H130").Value = 17857.143
$ws.Range("I130").Value = 5000
$ws.Range("J130").Value = 20000
$ws.Range("K130").Value = 15000
$ws.Range("L130").Value = 60000
$ws.Range("M130").Value = -9980
$ws.Range("N130").Value = -70040

$ws.Range("H131").Value = 2251.8
$ws.Range("J131").Value = 2819.6667
$ws.Range("L131").Value = 8459.000100000001
$ws.Range("N131").Value = -18539.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 775.6667
$ws.Range("I2").Value = 911.7368
$ws.Range("K2").Value = 911.7368
$ws.Range("M2").Value = -798.7368

$ws.Range("H97").Value = 11223
$ws.Range("I97").Value = 12829.167
$ws.Range("J97").Value = 5440.8
$ws.Range("K97").Value = 12829.167
$ws.Range("L97").Value = 5440.8
$ws.Range("M97").Value = -12333.167
$ws.Range("N97").Value = -6432.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 38770.8
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 1501.6
$ws.Range("I46").Value = 978.8
$ws.Range("J46").Value = 2024.4
$ws.Range("K46").Value = 978.8
$ws.Range("L46").Value = 2024.4
$ws.Range("M46").Value = -790.8
$ws.Range("N46").Value = -2400.4

$ws.Range("H93").Value = 5577.1724
$ws.Range("I93").Value = 6762.1816
$ws.Range("J93").Value = 1852.8572
$ws.Range("K93").Value = 6762.1816
$ws.Range("L93").Value = 1852.8572
$ws.Range("M93").Value = -5514.1816
$ws.Range("N93").Value = -4348.8572

$ws.Range("H122").Value = 6746.96
$ws.Range("I122").Value = 5037.5
$ws.Range("J122").Value = 11142.714
$ws.Range("K122").Value = 15112.5
$ws.Range("L122").Value = 33428.142
$ws.Range("M122").Value = -12662.5
$ws.Range("N122").Value = -38328.142

$ws.Range("H132").Value = 468448.44
$ws.Range("I132").Value = 828936.2
$ws.Range("K132").Value = 2486808.6
$ws.Range("M132").Value = -2484278.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 48333.332
$ws.Range("J99").Value = 48333.332
$ws.Range("L99").Value = 48333.332
$ws.Range("N99").Value = -54323.332

$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180

$ws.Range("H132").Value = 3924.7632
$ws.Range("I132").Value = 3706.4849
$ws.Range("J132").Value = 5365.4
$ws.Range("K132").Value = 11119.4547
$ws.Range("L132").Value = 16096.2
$ws.Range("M132").Value = -8589.4547
$ws.Range("N132").Value = -21156.2
